$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): right-answer marks and wrong-answer penalty updates
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): aggregated marks and penalty updates
$ws.Range("B12").Value = 65
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "59.0/140"
